$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.683.94'
$ws.Range('E2').Value = '  -1.52%  '
$ws.Range('D3').Value = '1.594.94'
$ws.Range('E3').Value = '  -1.79%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '211.48'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('E9').Value = '  -1.80%  '
$ws.Range('D10').Value = '19.69'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('D11').Value = '0.0836'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').Value = '1.818.03'
$ws.Range('E12').Value = '  -1.79%  '
$ws.Range('D13').Value = '1.628.94'
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  -3.24%  '
$ws.Range('D16').Value = '64.80'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').Value = '26.647.46'
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('D19').Value = '209.42'
$ws.Range('E19').Value = '  -1.92%  '
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').Value = '6.70'
$ws.Range('E21').Value = '  -2.19%  '
$ws.Range('E22').Value = '  -2.41%  '
$ws.Range('D23').Value = '2.33'
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('D24').Value = '8.89'
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('D25').Value = '146.77'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  -3.34%  '
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('D29').Value = '15.35'
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('D30').Value = '0.0504'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('E31').Value = '  -1.43%  '
$ws.Range('E32').Value = '  -3.70%  '
$ws.Range('D33').Value = '0.664'
$ws.Range('E33').Value = '  -9.26%  '
$ws.Range('D34').Value = '2.93'
$ws.Range('E34').Value = '  -2.57%  '
$ws.Range('D35').Value = '1.291.52'
$ws.Range('E35').Value = '  -5.02%  '
$ws.Range('D36').Value = '2.44'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('E37').Value = '  -5.33%  '
$ws.Range('E38').Value = '  -3.15%  '
$ws.Range('D39').Value = '0.836'
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '2.20'
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '5.36'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').Value = '63.57'
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('D45').Value = '1.730.75'
$ws.Range('D46').Value = '89.84'
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('D47').Value = '0.875'
$ws.Range('E47').Value = '  +0.77%  '
$ws.Range('E48').Value = '  -1.62%  '
$ws.Range('D49').Value = '0.0983'
$ws.Range('E49').Value = '  -2.61%  '
$ws.Range('D50').Value = '0.0504'
$ws.Range('E50').Value = '  -1.77%  '
$ws.Range('D51').Value = '7.47'
$ws.Range('E51').Value = '  -2.57%  '
